# Add new weekly ranking sheet "2026-02-04" at the end of the workbook,
# populate it with header + 100 ranked rows, and apply the existing
# "low volume" highlight style (yellow fill, same as other sheets) to
# the volume cells whose value is <= 3.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet as the last tab -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2026-02-04"

# --- 2. Header row -------------------------------------------------------
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "volume"
$ws.Range("D1").Value = "publisher"

# --- 3. Data rows (rank, title, volume) ----------------------------------
$rows = @(
    @(1,'スーパーの裏でヤニ吸うふたり 通常版',8),
    @(2,'メダリスト',14),
    @(3,'ブルーロック',37),
    @(4,'平和の国の島崎へ',11),
    @(5,'望郷太郎',14),
    @(6,'主人恋日記',12),
    @(7,'信者ゼロの女神サマと始める異世界攻略',1),
    @(8,'シャングリラ・フロンティア ~クソゲーハンター、神ゲーに挑まんとす~',25),
    @(9,'大正學生愛妻家',4),
    @(10,'フシノカミ 1 ~辺境から始める文明再生記~',1),
    @(11,'おはらい箱の天才付与術師は、辺境で悠々自適に暮らしたい~万能付与術で気付いたら辺境が世界最強の快適拠点になっていた~',1),
    @(12,'ひとりぼっちの異世界攻略',26),
    @(13,'煙たい話',2),
    @(14,'ダンダダン',22),
    @(15,'「壇蜜」',2),
    @(16,'高度に発達した医学は魔法と区別がつかない',11),
    @(17,'青野くんに触りたいから死にたい',14),
    @(18,'陸上自衛隊特務諜報機関 別班の犬',9),
    @(19,'Mr.CB',17),
    @(20,'僕の心のヤバイやつ',13),
    @(21,'葬送のフリーレン',15),
    @(22,'ラガーにゃん 1~猫でもわかるラグビー入門~',1),
    @(23,'最期の願いに月が泣く',1),
    @(24,'ひとりぼっちの異世界攻略',1),
    @(25,'フェルマーの料理',7),
    @(26,'貴族転生 ~恵まれた生まれから最強の力を得る~',10),
    @(27,'TSUYOSHI 誰も勝てない、アイツには',29),
    @(28,'ONE PIECE',113),
    @(29,'#神奈川に住んでるエルフ',1),
    @(30,'うららか',1),
    @(31,'多聞くん今どっち!?',13),
    @(32,'はじめの一歩',145),
    @(33,'るろうに剣心―明治剣客浪漫譚・北海道編―',10),
    @(34,'弱虫ペダル',98),
    @(35,'呪術廻戦≡(モジュロ)',1),
    @(36,'追放されたパシリ、買い物スキルSSSで装備無双 ~買ったモノを超強化して最強パーティー目指します~第1話',1),
    @(37,'婚約者が浮気しているようなんですけど私は流行りの悪役令嬢ってことであってますか?',3),
    @(38,'ひびワれこうじつ',1),
    @(39,'俺、悪役騎士団長に転生する。',1),
    @(40,'俺、悪役騎士団長に転生する。',2),
    @(41,'俺、悪役騎士団長に転生する。',3),
    @(42,'銀河特急 ミルキー☆サブウェイ',1),
    @(43,'最期の願いに月が泣く',2),
    @(44,'最期の願いに月が泣く',3),
    @(45,'黒騎士様から全力で溺愛されていますが、すごもり聖女は今日も引きこもりたい!',1),
    @(46,'バラバラ夫婦~手足をなくした夫はまだ生きてる1',1),
    @(47,'サラリーマンが異世界に行ったら四天王になった話',1),
    @(48,'異世界整体師 ~美女も亜人も魔物も竜も、お前ら全員揉みほぐす!!~',2),
    @(49,'限界女神 異世界転生窓口はアットホームな職場です',2),
    @(50,'追放された冒険者たちはHなダンジョンでなりあがる!1',1),
    @(51,'転生貴族、鑑定スキルで成り上がる ~弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた~',20),
    @(52,'ダーウィン事変',10),
    @(53,'ザ・ファブル The third secret',3),
    @(54,'二階堂地獄ゴルフ',11),
    @(55,'追放貴族は、外れスキルで英霊たちと辺境領地を再興する~英霊たちを召喚したら慕われたので、最強領地を作り上げます~',7),
    @(56,'Shrink~精神科医ヨワイ~',17),
    @(57,'彼女、お借りします',44),
    @(58,'ワンパンマン',34),
    @(59,'ワンパンマン',35),
    @(60,'スーパーの裏でヤニ吸うふたり',2),
    @(61,'ハズレスキルが覚醒したので、虐げられた過去を自由な人生に反転(リバース)します!',1),
    @(62,'次に買うマンガ、この1話で決めよう! ~ハマったら抜け出せない! 王道ファンタジー編~',5),
    @(63,'冗談で女勇者たちを口説いた大魔王、攫われて強制新婚生活1',1),
    @(64,'パラダイスシフト1',1),
    @(65,'黒髪清楚の冷酷美少女を助けたら、俺と二人きりの時だけデレるようになった件1',1),
    @(66,'この結婚は誰がために',3),
    @(67,'追放された冒険者たちはHなダンジョンでなりあがる!2',2),
    @(68,'田舎のホームセンター男の自由な異世界生活',15),
    @(69,'シャドーハウス',21),
    @(70,'邪神の弁当屋さん',4),
    @(71,'空挺ドラゴンズ',21),
    @(72,'スーパーの裏でヤニ吸うふたり',4),
    @(73,'スーパーの裏でヤニ吸うふたり 通常版',6),
    @(74,'スーパーの裏でヤニ吸うふたり',7),
    @(75,'貴族転生 ~恵まれた生まれから最強の力を得る~',9),
    @(76,'薬屋のひとりごと',16),
    @(77,'スーパーの裏でヤニ吸うふたり 通常版',3),
    @(78,'スーパーの裏でヤニ吸うふたり',5),
    @(79,'ポイントギフター《経験値分配能力者》の異世界最強ソロライフ ~ブラックギルドから解放された男は万能最強職として無双する~(コミック)',7),
    @(80,'学園騎士のレベルアップ!レベル1000超えの転生者、落ちこぼれクラスに入学。そして、(コミック)',12),
    @(81,'フシノカミ 10 ~辺境から始める文明再生記~',10),
    @(82,'葬送のフリーレン',14),
    @(83,'金色のガッシュ!!',6),
    @(84,'SAKAMOTO DAYS',25),
    @(85,'カグラバチ',10),
    @(86,'乙女椿は笑わない',13),
    @(87,'呪術廻戦',27),
    @(88,'呪術廻戦',28),
    @(89,'呪術廻戦',29),
    @(90,'呪術廻戦',30),
    @(91,'呪術廻戦',26),
    @(92,'ビジネス婚ー好きになったら離婚しますー',1),
    @(93,'追加ダメージ1しか与えられない底辺配信者、双剣に持ち替えた途端に無双する1',1),
    @(94,'僕らの恋は夢うつつ',1),
    @(95,'兄様のライバルに目をつけられたようです(エンジェライトコミックス)',1),
    @(96,'恋フレ ~恋人未満がちょうどいい~',1),
    @(97,'ひびワれこうじつ',2),
    @(98,'ひびワれこうじつ',3),
    @(99,'銀河特急 ミルキー☆サブウェイ',2),
    @(100,'銀河特急 ミルキー☆サブウェイ',3)
)

$data = New-Object 'object[,]' $rows.Count,3
for ($i = 0; $i -lt $rows.Count; $i++) {
    $data[$i,0] = $rows[$i][0]
    $data[$i,1] = $rows[$i][1]
    $data[$i,2] = $rows[$i][2]
}

$lastRow = 1 + $rows.Count
$ws.Range("A2:C$lastRow").Value = $data

# --- 4. Highlight volume cells whose value is <= 3 -----------------------
# Reuse the exact fill style already used elsewhere in the workbook for
# this convention by copying an existing highlighted cell's format.
$styleSource = $wb.Worksheets.Item("2026-01-28").Range("C15")
$styleSource.Copy()

$highlightCells = @("C8","C11","C12","C14","C16","C23","C24","C25","C30","C31","C36","C37","C38","C39","C40","C41","C42","C43","C44","C45","C46","C47","C48","C49","C50","C51","C54","C61","C62","C64","C65","C66","C67","C68","C78","C93","C94","C95","C96","C97","C98","C99","C100","C101")
foreach ($addr in $highlightCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

Write-Output "Added sheet 2026-02-04 with $($rows.Count) rows; highlighted $($highlightCells.Count) low-volume cells."
